$d = $word.ActiveDocument

function Replace-ParagraphText($para, [string]$newText) {
    # A direct `Range.Text = ...` assignment in this host only overwrites the
    # first run of a range that spans multiple runs, so route existing
    # (non-empty) paragraph text replacement through Find/Replace instead,
    # which correctly rewrites text spanning multiple runs. The paragraph
    # Range.Text carries a trailing paragraph-mark (\r) that Find cannot
    # match literally, so strip it first.
    $rng = $para.Range
    $oldText = $rng.Text
    $oldText = $oldText.Substring(0, $oldText.Length - 1)
    if ($oldText.Length -eq 0) {
        # Brand-new empty paragraph (single empty run) -- a direct Range.Text
        # assignment works fine here (no multi-run merge issue), and Find
        # cannot search for an empty string anyway.
        $rng.Text = $newText
        return
    }
    $ok = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Replace-ParagraphText failed to find [$oldText]"
    }
}

function Assert-Text($para, [string]$expected) {
    $actual = $para.Range.Text
    $actual = $actual.Substring(0, $actual.Length - 1)
    if ($actual -ne $expected) {
        throw "Assertion failed: expected [$expected] but found [$actual]"
    }
}

# --- Title paragraph ---
# "Instructions (more suggestions – rushed and not tested) "
#   -> "Data Rename tool installation " + "Instructions "
$p1 = $d.Paragraphs.Item(1)
Assert-Text $p1 "Instructions (more suggestions – rushed and not tested) "
Replace-ParagraphText $p1 "Data Rename tool installation Instructions "

# --- Second paragraph ---
# "Double click the add-in to install.  " -> "Beta version January 2014."
$p2 = $d.Paragraphs.Item(2)
Assert-Text $p2 "Double click the add-in to install.  "
Replace-ParagraphText $p2 "Beta version January 2014."

# --- Third paragraph is the blank spacer paragraph; left untouched. ---

# --- First ListParagraph bullet ---
# "Double click the add-in to install.  " -> "Remove all" + " previous versions"
$p4 = $d.Paragraphs.Item(4)
Assert-Text $p4 "Double click the add-in to install.  "
Replace-ParagraphText $p4 "Remove all previous versions"

# --- Second ListParagraph bullet ---
# "Copy the lookup database here: C:\...\Resources" -> "Double click the add-in to install.  "
$p5 = $d.Paragraphs.Item(5)
Assert-Text $p5 "Copy the lookup database here: C:\Program Files (x86)\ArcGIS\Desktop10.1\bin\Resources"
Replace-ParagraphText $p5 "Double click the add-in to install.  "

# --- Third ListParagraph bullet (currently has no numPr) ---
# Old: "*You may have to manually create the 'Resources' folder above"
# New: a numbered paragraph "Create a new directory in crash the crash move
#      folder ..GIS\2_Active_Data\200_data_name_lookup" followed by a plain
#      ListParagraph (no numPr) "and copy the lookup csv files into it".
$p6 = $d.Paragraphs.Item(6)
Assert-Text $p6 "*You may have to manually create the ‘Resources’ folder above"

# Insert the trailing "and copy..." paragraph right after p6 first, so it
# inherits p6's current pPr shape (ListParagraph, no numPr) -- matching the
# desired final formatting for that trailing paragraph.
$p6.Range.InsertParagraphAfter()
$pCopy = $d.Paragraphs.Item(7)
Replace-ParagraphText $pCopy "and copy the lookup csv files into it"

# Insert the new numbered paragraph right after paragraph 5 (which already has
# the ListParagraph + numPr pPr), so the new paragraph inherits that numPr.
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertParagraphAfter()
$pCreate = $d.Paragraphs.Item(6)
Replace-ParagraphText $pCreate "Create a new directory in crash the crash move folder ..GIS\2_Active_Data\200_data_name_lookup"

# Delete the now-redundant old paragraph ("*You may have to manually create
# the 'Resources' folder above"), which has shifted down to index 7 after the
# two insertions above (6: new "Create a new directory..." paragraph, 7: old
# paragraph, 8: new "and copy..." paragraph).
$pOld = $d.Paragraphs.Item(7)
Assert-Text $pOld "*You may have to manually create the ‘Resources’ folder above"
$pOld.Range.Delete()

Write-Output "Edit complete. Paragraph count: $($d.Paragraphs.Count)"
